$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) unaffected in content, only P1:R1 share-string indices shift
# (values stay the same text, just shared string table was deduped/reordered) ---
$ws.Range("P1").Value = "Volume"
$ws.Range("Q1").Value = "Fragment Size"
$ws.Range("R1").Value = "Read Length"

# --- Row 2: replace unique/test-only values, add styles to C2/D2 ---
$ws.Range("A2").Value = "(unique value filled in by the test)"
$ws.Range("B2").Value = "(unique value filled in by the test)"
$ws.Range("C2").Value = "(unique value filled in by the test)"
$ws.Range("D2").Value = "(same as broad sample)"
$ws.Range("C2").Style = $ws.Range("A2").Style
$ws.Range("D2").Style = $ws.Range("A2").Style

$ws.Range("J2").Value = "COLB-123"
$ws.Range("K2").Value = "COLAB-P-234"
$ws.Range("L2").Value = "BP-ID-567"
$ws.Range("M2").Value = "M"
$ws.Range("N2").Value = "Canine"
$ws.Range("O2").Value = "lsid:1"

# --- Row 3: brand-new data row replacing the old blank formatting-only row ---
$ws.Range("A3").Value = "(unique value filled in by the test)"
$ws.Range("B3").Value = "(unique value filled in by the test)"
$ws.Range("C3").Value = "(unique value filled in by the test)"
$ws.Range("D3").Value = "(unique value filled in by the test)"
$ws.Range("A3").Style = $ws.Range("A2").Style
$ws.Range("B3").Style = $ws.Range("A2").Style
$ws.Range("C3").Style = $ws.Range("A2").Style
$ws.Range("D3").Style = $ws.Range("A2").Style

$ws.Range("E3").Value = "Illumina_P5-Nijow_P7-Waren"
$ws.Range("H3").Value = "DEV-6796"
$ws.Range("I3").Value = "DEV-6815, DEV-6816"
$ws.Range("J3").Value = "COLB-124"
$ws.Range("K3").Value = "COLAB-P-235"
$ws.Range("L3").Value = "BP-ID-568"
$ws.Range("M3").Value = "F"
$ws.Range("N3").Value = "Feline"
$ws.Range("O3").Value = "lsid:2"
$ws.Range("P3").Value = 62
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 4

# Row 3 previously had a taller custom row height and a lone formatted (wrap-text)
# placeholder cell at I3; the new row uses normal height/default formatting.
$ws.Rows.Item(3).RowHeight = $ws.Rows.Item(2).RowHeight

# --- View / selection tweaks ---
$ws.Range("G3").Select() | Out-Null
